$p = $ppt.ActivePresentation

# --- 1. Fix the DSW label on the map (was mislabeled "MRN") ---
$slide = $p.Slides.Item(1)
for ($k = 1; $k -le $slide.Shapes.Count; $k++) {
    $shp = $slide.Shapes.Item($k)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq "MRN") {
            $shp.TextFrame.TextRange.Text = "DSW"
        }
    }
}

# --- 2. Roll the "last updated" date shown on every layout/master forward a day ---
$oldDate = "6/21/19"
$newDate = "6/22/19"

function Update-DatePlaceholder($container) {
    for ($j = 1; $j -le $container.Shapes.Count; $j++) {
        $shp = $container.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText -and $shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    Update-DatePlaceholder $master.CustomLayouts.Item($i)
}

Write-Output "Applied DSW label fix and date refresh."
